# Update "想去人数" (want-to-go count) values in column F across all sheets,
# reflecting newly generated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1772
$ws.Range("F4").Value  = 447
$ws.Range("F7").Value  = 626
$ws.Range("F8").Value  = 322
$ws.Range("F9").Value  = 1711
$ws.Range("F10").Value = 356
$ws.Range("F12").Value = 804
$ws.Range("F13").Value = 329
$ws.Range("F15").Value = 12724
$ws.Range("F16").Value = 12746
$ws.Range("F17").Value = 946
$ws.Range("F18").Value = 740
$ws.Range("F21").Value = 51
$ws.Range("F22").Value = 538
$ws.Range("F23").Value = 1993
$ws.Range("F24").Value = 25
$ws.Range("F27").Value = 668

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 4
$ws.Range("F8").Value  = 127
$ws.Range("F10").Value = 74

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 84
$ws.Range("F3").Value = 162

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 84
$ws.Range("F3").Value  = 162
$ws.Range("F5").Value  = 1772
$ws.Range("F6").Value  = 447
$ws.Range("F7").Value  = 4
$ws.Range("F11").Value = 626
$ws.Range("F13").Value = 322
$ws.Range("F14").Value = 1711
$ws.Range("F15").Value = 356
$ws.Range("F17").Value = 804
$ws.Range("F18").Value = 329
$ws.Range("F21").Value = 12724
$ws.Range("F22").Value = 12746
$ws.Range("F23").Value = 946
$ws.Range("F24").Value = 740
$ws.Range("F27").Value = 51
$ws.Range("F28").Value = 538
$ws.Range("F31").Value = 1993
$ws.Range("F32").Value = 25
$ws.Range("F33").Value = 127
$ws.Range("F37").Value = 668
$ws.Range("F38").Value = 74
